$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the day label (column A) for every row of each 4-row block ---
# (previously only the first row of each block carried the label; now every row does)
$days = @("Pon.", "Wt.", "Śr.", "Czw.", "Pt.", "Sob.", "Niedz.")
$row = 2
foreach ($d in $days) {
    for ($i = 0; $i -lt 4; $i++) {
        $ws.Cells.Item($row + $i, 1).Value = $d
    }
    # Copy the border formatting of the block's first row onto the other three
    # rows of the block (Left+Right+Top thin border), reusing the existing style.
    $ws.Range("A" + $row).Copy()
    $ws.Range("A" + ($row + 1) + ":A" + ($row + 3)).PasteSpecial(-4122)
    $row += 4
}

# Row 29 closes the column-A box, so it gets the fully-boxed (all 4 sides) border,
# matching the header row's style.
$ws.Range("A1").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Niedz."

# --- Remove the empty "Komentarze" column (I) ---
$ws.Columns.Item(9).Delete()

# --- Restore the view: scrolled back to column A, selection on I1 ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I1").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

